$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several quarterly balance sheets were concatenated into this single sheet;
# rows 64 ("Perdas pela Nao Recuperabilidade de Ativos") and 79 ("Part. de
# Acionistas Nao Controladores") only had placeholder zeros for the periods
# where the line item did not exist in the source statement. Those zeros are
# replaced with genuinely blank/empty cells (matching the already-empty B and
# D columns) so the concatenated sheet doesn't imply real zero values.
$rows = @(64, 79)

foreach ($r in $rows) {
    $single = $ws.Range("C$r")
    $single.Value = "'"
    $single.ClearFormats()

    $rest = $ws.Range("E$r`:AG$r")
    $rest.Value = "'"
    $rest.ClearFormats()
}
